$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percentage-formatted text cells: force text storage so Excel does not
# auto-convert "79%" into a numeric 0.79 value with a percent number format.
$percentCells = @("H3", "H4", "H5", "H6", "H8", "H9", "H12", "H20", "H22", "H26", "H29", "H33")
foreach ($cellRef in $percentCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("H3").Value = '79%'
$ws.Range("H4").Value = '66%'
$ws.Range("H5").Value = '70%'
$ws.Range("H6").Value = '70%'
$ws.Range("H8").Value = '84%'
$ws.Range("H9").Value = '99%'
$ws.Range("H12").Value = '85%'
$ws.Range("H20").Value = '72%'
$ws.Range("H22").Value = '86%'
$ws.Range("H26").Value = '77%'
$ws.Range("H29").Value = '78%'
$ws.Range("H33").Value = '84%'

# Plain text / date-time / measurement cells (kept as text automatically).
$ws.Range("E2").Value = '2026-02-05 21:18:03'
$ws.Range("E3").Value = '2026-02-05 21:18:06'
$ws.Range("E4").Value = '2026-02-05 21:18:08'
$ws.Range("J4").Value = '989.6 hPa'
$ws.Range("E5").Value = '2026-02-05 21:18:11'
$ws.Range("O5").Value = '10.1 °C'
$ws.Range("E6").Value = '2026-02-05 21:18:14'
$ws.Range("O6").Value = '13.1 °C'
$ws.Range("E7").Value = '2026-02-05 21:18:16'
$ws.Range("E8").Value = '2026-02-05 21:18:19'
$ws.Range("E9").Value = '2026-02-05 21:18:22'
$ws.Range("O9").Value = '2.4 °C'
$ws.Range("E10").Value = '2026-02-05 21:18:25'
$ws.Range("E11").Value = '2026-02-05 21:18:27'
$ws.Range("E12").Value = '2026-02-05 21:18:30'
$ws.Range("O12").Value = '10.4 °C'
$ws.Range("E13").Value = '2026-02-05 21:18:33'
$ws.Range("E14").Value = '2026-02-05 21:18:35'
$ws.Range("I14").Value = '7.9 mm'
$ws.Range("E15").Value = '2026-02-05 21:18:38'
$ws.Range("L15").Value = '52.9 km/h - 261º 20:35 TU'
$ws.Range("O15").Value = '8.7 °C'
$ws.Range("E16").Value = '2026-02-05 21:18:41'
$ws.Range("E17").Value = '2026-02-05 21:18:44'
$ws.Range("I17").Value = '8.6 mm'
$ws.Range("E18").Value = '2026-02-05 21:18:47'
$ws.Range("E19").Value = '2026-02-05 21:18:49'
$ws.Range("E20").Value = '2026-02-05 21:18:52'
$ws.Range("E21").Value = '2026-02-05 21:18:55'
$ws.Range("E22").Value = '2026-02-05 21:18:58'
$ws.Range("O22").Value = '9.1 °C'
$ws.Range("E23").Value = '2026-02-05 21:19:01'
$ws.Range("O23").Value = '8.2 °C'
$ws.Range("E24").Value = '2026-02-05 21:19:03'
$ws.Range("E25").Value = '2026-02-05 21:19:06'
$ws.Range("J25").Value = '993.8 hPa'
$ws.Range("O25").Value = '0.8 °C'
$ws.Range("E26").Value = '2026-02-05 21:19:09'
$ws.Range("E27").Value = '2026-02-05 21:19:11'
$ws.Range("J27").Value = '990.0 hPa'
$ws.Range("E28").Value = '2026-02-05 21:19:14'
$ws.Range("O28").Value = '2.8 °C'
$ws.Range("E29").Value = '2026-02-05 21:19:16'
$ws.Range("O29").Value = '9.5 °C'
$ws.Range("E30").Value = '2026-02-05 21:19:19'
$ws.Range("M30").Value = '0.5 °C 19:34 TU'
$ws.Range("O30").Value = '-1.6 °C'
$ws.Range("E31").Value = '2026-02-05 21:19:22'
$ws.Range("I31").Value = '19.5 mm'
$ws.Range("E32").Value = '2026-02-05 21:19:24'
$ws.Range("O32").Value = '12.1 °C'
$ws.Range("E33").Value = '2026-02-05 21:19:27'
$ws.Range("O33").Value = '9.6 °C'
$ws.Range("E34").Value = '2026-02-05 21:19:30'
$ws.Range("O34").Value = '4.4 °C'
$ws.Range("E35").Value = '2026-02-05 21:19:33'
$ws.Range("I35").Value = '5.5 mm'
$ws.Range("E36").Value = '2026-02-05 21:19:36'
$ws.Range("J36").Value = '992.5 hPa'
$ws.Range("O36").Value = '10.0 °C'
